$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting existing rows (3..44) down to (4..45)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new record.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(3, 4).Value = 45092
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat
$ws.Cells.Item(3, 5).Value = 15
$ws.Cells.Item(3, 6).Value = 100112044
$ws.Cells.Item(3, 7).Value = "Perejil"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 300
$ws.Cells.Item(3, 11).Value = 1800
$ws.Cells.Item(3, 12).Value = 2000
$ws.Cells.Item(3, 13).Value = 1900
$ws.Cells.Item(3, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(3, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(3, 16).Value = 950
$ws.Cells.Item(3, 17).Value = 2
$ws.Cells.Item(3, 18).Value = "Hortaliza"
